# Tweaks to word template doc
#
# This edits the styles used by the qsc_docx Word template:
#   - Normal: paragraphs become fully justified (w:jc val="both")
#   - Heading1/2/3: paragraph spacing before/after is tightened up
#   - TableCaption: no longer needs its own explicit justification
#     (it now inherits "both" from Normal)
#   - the abstractNum "w:tmpl" fingerprints in numbering.xml are
#     refreshed, as Word does whenever it resaves the numbering part
#
# The rsid stamps on the touched styles are bumped too, mirroring what
# Word records for the edit session that made these changes.

$d = $word.ActiveDocument

function Replace-Once {
    param(
        [string]$Xml,
        [string]$Old,
        [string]$New,
        [string]$Label
    )
    if (-not $Xml.Contains($Old)) {
        throw ("Replace-Once: pattern not found for " + $Label)
    }
    return $Xml.Replace($Old, $New)
}

$xml = $d.WordOpenXML

# --- styles.xml -----------------------------------------------------

# Normal: new rsid + justify all paragraphs by default
$xml = Replace-Once $xml `
    '<w:style w:type="paragraph" w:default="1" w:styleId="Normal"><w:name w:val="Normal"/><w:qFormat/><w:rsid w:val="00E23021"/><w:rPr>' `
    '<w:style w:type="paragraph" w:default="1" w:styleId="Normal"><w:name w:val="Normal"/><w:qFormat/><w:rsid w:val="004C24F2"/><w:pPr><w:jc w:val="both"/></w:pPr><w:rPr>' `
    "Normal"

# Heading1: new rsid + spacing after=240 added
$xml = Replace-Once $xml `
    '<w:link w:val="Heading1Char"/><w:uiPriority w:val="9"/><w:qFormat/><w:rsid w:val="00E23021"/><w:pPr><w:keepNext/><w:keepLines/><w:spacing w:before="240"/><w:outlineLvl w:val="0"/></w:pPr>' `
    '<w:link w:val="Heading1Char"/><w:uiPriority w:val="9"/><w:qFormat/><w:rsid w:val="004C2F38"/><w:pPr><w:keepNext/><w:keepLines/><w:spacing w:before="240" w:after="240"/><w:outlineLvl w:val="0"/></w:pPr>' `
    "Heading1"

# Heading2: new rsid + spacing before/after=240
$xml = Replace-Once $xml `
    '<w:link w:val="Heading2Char"/><w:uiPriority w:val="9"/><w:semiHidden/><w:unhideWhenUsed/><w:qFormat/><w:rsid w:val="00E8211A"/><w:pPr><w:keepNext/><w:keepLines/><w:spacing w:before="40"/><w:outlineLvl w:val="1"/></w:pPr>' `
    '<w:link w:val="Heading2Char"/><w:uiPriority w:val="9"/><w:semiHidden/><w:unhideWhenUsed/><w:qFormat/><w:rsid w:val="004C2F38"/><w:pPr><w:keepNext/><w:keepLines/><w:spacing w:before="240" w:after="240"/><w:outlineLvl w:val="1"/></w:pPr>' `
    "Heading2"

# Heading3: new rsid + spacing before/after=120
$xml = Replace-Once $xml `
    '<w:link w:val="Heading3Char"/><w:uiPriority w:val="9"/><w:unhideWhenUsed/><w:qFormat/><w:rsid w:val="00E8211A"/><w:pPr><w:keepNext/><w:keepLines/><w:spacing w:before="40"/><w:outlineLvl w:val="2"/></w:pPr>' `
    '<w:link w:val="Heading3Char"/><w:uiPriority w:val="9"/><w:unhideWhenUsed/><w:qFormat/><w:rsid w:val="004C2F38"/><w:pPr><w:keepNext/><w:keepLines/><w:spacing w:before="120" w:after="120"/><w:outlineLvl w:val="2"/></w:pPr>' `
    "Heading3"

# Heading1Char: new rsid
$xml = Replace-Once $xml `
    '<w:link w:val="Heading1"/><w:uiPriority w:val="9"/><w:rsid w:val="00E23021"/><w:rPr>' `
    '<w:link w:val="Heading1"/><w:uiPriority w:val="9"/><w:rsid w:val="004C2F38"/><w:rPr>' `
    "Heading1Char"

# Heading2Char: new rsid
$xml = Replace-Once $xml `
    '<w:link w:val="Heading2"/><w:uiPriority w:val="9"/><w:semiHidden/><w:rsid w:val="00E8211A"/><w:rPr>' `
    '<w:link w:val="Heading2"/><w:uiPriority w:val="9"/><w:semiHidden/><w:rsid w:val="004C2F38"/><w:rPr>' `
    "Heading2Char"

# TableCaption: drop the now-redundant explicit justification (Normal covers it)
$xml = Replace-Once $xml `
    '<w:basedOn w:val="Normal"/><w:qFormat/><w:rsid w:val="00575994"/><w:pPr><w:jc w:val="both"/></w:pPr><w:rPr><w:iCs/>' `
    '<w:basedOn w:val="Normal"/><w:qFormat/><w:rsid w:val="00575994"/><w:rPr><w:iCs/>' `
    "TableCaption"

# Heading3Char: new rsid
$xml = Replace-Once $xml `
    '<w:link w:val="Heading3"/><w:uiPriority w:val="9"/><w:rsid w:val="00E8211A"/><w:rPr>' `
    '<w:link w:val="Heading3"/><w:uiPriority w:val="9"/><w:rsid w:val="004C2F38"/><w:rPr>' `
    "Heading3Char"

# --- numbering.xml ----------------------------------------------------
# Refresh the abstractNum w:tmpl fingerprints (cosmetic internal ids Word
# regenerates whenever it resaves the numbering part).

$tmplMap = @{
    "A3C2CCD6" = "A18CE72C";
    "ACA850D6" = "FE6E5ECC";
    "9394050A" = "2728ABDC";
    "8DF204AE" = "A2CE2A26";
    "EDB626B6" = "BF965796";
    "9C80746E" = "929C17EA";
    "441A2F00" = "D2D61544";
    "88209878" = "2478778A";
    "CB46F3C2" = "581C9098";
    "9FF04790" = "8418F24A";
}

foreach ($old in $tmplMap.Keys) {
    $new = $tmplMap[$old]
    $xml = Replace-Once $xml ('<w:tmpl w:val="' + $old + '"/>') ('<w:tmpl w:val="' + $new + '"/>') ("tmpl " + $old)
}

$d.WordOpenXML = $xml

Write-Output "applied qsc_docx template tweaks"
